$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = "CheckBox (BookApp)"
$ws.Range("C27").Value = "RadioButton(Traffic Signal App)"
$ws.Range("C28").Value = "progressBar View"

# Append new strings in the order needed to match target shared string indices 32-38
$ws.Range("C29").Value = "Basic View 2"
$ws.Range("C30").Value = "AutoCompleteTextView"
$ws.Range("C31").Value = "UsingPickerView"
$ws.Range("C34").Value = "ListView"
$ws.Range("C35").Value = "SpinnerView"
$ws.Range("C32").Value = "DatePicker"
$ws.Range("C33").Value = "TimePicker"

$ws.Range("D28").Value = "LabBook"
$ws.Range("D29").Value = "Journal"
$ws.Range("D30").Value = "LabBook"
$ws.Range("D31").Value = "LabBook"
$ws.Range("D32").Value = "Journal"
$ws.Range("D33").Value = "LabBook"
$ws.Range("D34").Value = "Journal"
$ws.Range("D35").Value = "Journal"

$ws.Range("A28").ClearContents()

# Extend the conditional formatting range to cover the newly added rows
$fcs = $ws.Range("D3:D28").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D3:D35"))

# Scroll/selection bookkeeping to match the saved view state
$ws.Range("C32").Select() | Out-Null
